$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "76529250-6/0"
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 16956.9559
$ws.Range("F2").Value = "18/03/2025"
$ws.Range("G2").Value = "18/03/2025"
$ws.Range("H2").Value = 254354

# Row 3
$ws.Range("A3").Value = "96586750-3/0"
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 16956.9559
$ws.Range("F3").Value = "18/03/2025"
$ws.Range("G3").Value = "18/03/2025"
$ws.Range("H3").Value = 135656
